$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6665.6665
$ws.Range("J64").Value = 7499.75
$ws.Range("L64").Value = 7499.75
$ws.Range("N64").Value = -7995.75
$ws.Range("H67").Value = 6665.6665
$ws.Range("J67").Value = 7499.75
$ws.Range("L67").Value = 7499.75
$ws.Range("N67").Value = -9215.75
$ws.Range("H76").Value = 3102
$ws.Range("I76").Value = 3069.6667
$ws.Range("J76").Value = 3199
$ws.Range("K76").Value = 3069.6667
$ws.Range("L76").Value = 3199
$ws.Range("M76").Value = -2754.6667
$ws.Range("N76").Value = -3829
$ws.Range("H79").Value = 3102
$ws.Range("I79").Value = 3069.6667
$ws.Range("J79").Value = 3199
$ws.Range("K79").Value = 3069.6667
$ws.Range("L79").Value = 3199
$ws.Range("M79").Value = -1977.6667
$ws.Range("N79").Value = -5383
$ws.Range("H106").Value = 18798.75
$ws.Range("I106").Value = 15065
$ws.Range("K106").Value = 15065
$ws.Range("M106").Value = -14434
$ws.Range("H137").Value = 2171.2
$ws.Range("I137").Value = 2016
$ws.Range("J137").Value = 2248.8
$ws.Range("K137").Value = 6048
$ws.Range("L137").Value = 6746.400000000001
$ws.Range("M137").Value = -3498
$ws.Range("N137").Value = -11846.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()  # was -659
$ws.Range("N33").ClearContents()  # was -25558
$ws.Range("H74").Value = 2808.4546
$ws.Range("I74").Value = 1831.6666
$ws.Range("J74").Value = 3980.6
$ws.Range("K74").Value = 1831.6666
$ws.Range("L74").Value = 3980.6
$ws.Range("M74").Value = -957.6666
$ws.Range("N74").Value = -5728.6
$ws.Range("H77").Value = 2808.4546
$ws.Range("I77").Value = 1831.6666
$ws.Range("J77").Value = 3980.6
$ws.Range("K77").Value = 9158.333000000001
$ws.Range("L77").Value = 19903
$ws.Range("M77").Value = -4790.333000000001
$ws.Range("N77").Value = -28639
$ws.Range("H88").Value = 3109.6
$ws.Range("J88").Value = 3109.6
$ws.Range("L88").Value = 3109.6
$ws.Range("N88").Value = -3921.6
$ws.Range("H91").Value = 3109.6
$ws.Range("J91").Value = 3109.6
$ws.Range("L91").Value = 3109.6
$ws.Range("N91").Value = -5917.6
$ws.Range("H122").Value = 1975
$ws.Range("I122").Value = 1982.1428
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 5946.428400000001
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -3496.428400000001
$ws.Range("N122").Value = -10750
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 392.27274
$ws.Range("J7").Value = 155.25
$ws.Range("L7").Value = 155.25
$ws.Range("N7").Value = -381.25
$ws.Range("H31").Value = 2519.125
$ws.Range("I31").Value = 1525.5
$ws.Range("J31").Value = 5500
$ws.Range("K31").Value = 1525.5
$ws.Range("L31").Value = 5500
$ws.Range("M31").Value = -1230.5
$ws.Range("N31").Value = -6090
$ws.Range("H34").Value = 2519.125
$ws.Range("I34").Value = 1525.5
$ws.Range("J34").Value = 5500
$ws.Range("K34").Value = 1525.5
$ws.Range("L34").Value = 5500
$ws.Range("M34").Value = -1323.5
$ws.Range("N34").Value = -5904
$ws.Range("H39").Value = 9999.666999999999
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 12499.5
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 12499.5
$ws.Range("M39").Value = -4609
$ws.Range("N39").Value = -13281.5
$ws.Range("H49").Value = 9999.666999999999
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 12499.5
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 12499.5
$ws.Range("M49").Value = -4818
$ws.Range("N49").Value = -12863.5
$ws.Range("H86").Value = 4899
$ws.Range("I86").Value = 4848.5
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4848.5
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3725.5
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4899
$ws.Range("I89").Value = 4848.5
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 24242.5
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -18626.5
$ws.Range("N89").Value = -36232
$ws.Range("H132").Value = 2462.3684
$ws.Range("I132").Value = 2236.75
$ws.Range("K132").Value = 6710.25
$ws.Range("M132").Value = -4180.25
$ws.Range("H134").Value = 3136.5
$ws.Range("I134").Value = 3136.5
$ws.Range("K134").Value = 9409.5
$ws.Range("M134").Value = -6874.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 127.53846
$ws.Range("I12").Value = 167.25
$ws.Range("K12").Value = 501.75
$ws.Range("M12").Value = -328.75
$ws.Range("H26").Value = 796.2917
$ws.Range("I26").Value = 441.8421
$ws.Range("J26").Value = 2143.2
$ws.Range("K26").Value = 1325.5263
$ws.Range("L26").Value = 6429.599999999999
$ws.Range("M26").Value = -1037.5263
$ws.Range("N26").Value = -7005.599999999999
$ws.Range("H56").Value = 5413.8247
$ws.Range("I56").Value = 5413.8247
$ws.Range("K56").Value = 5413.8247
$ws.Range("M56").Value = -4883.8247
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 175.35715
$ws.Range("J2").Value = 76.2
$ws.Range("L2").Value = 76.2
$ws.Range("N2").Value = -302.2
$ws.Range("H102").Value = 3236.2856
$ws.Range("I102").Value = 3331.3845
$ws.Range("K102").Value = 3331.3845
$ws.Range("M102").Value = -1709.3845
$ws.Range("H122").Value = 3889
$ws.Range("I122").Value = 3098.5715
$ws.Range("J122").Value = 5733.3335
$ws.Range("K122").Value = 9295.7145
$ws.Range("L122").Value = 17200.0005
$ws.Range("M122").Value = -6845.7145
$ws.Range("N122").Value = -22100.0005
$ws.Range("H123").Value = 8944
$ws.Range("J123").Value = 8944
$ws.Range("L123").Value = 8944
$ws.Range("N123").Value = -13844
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 514.8333
$ws.Range("I9").Value = 497.8
$ws.Range("K9").Value = 497.8
$ws.Range("M9").Value = -273.8
$ws.Range("H40").Value = 6279.1
$ws.Range("I40").Value = 6279.1
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6279.1
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6143.1
$ws.Range("N40").ClearContents()  # was -6072
$ws.Range("H46").Value = 1543.25
$ws.Range("I46").Value = 1543.25
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1543.25
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1355.25
$ws.Range("N46").ClearContents()  # was -3776
$ws.Range("H55").Value = 1460
$ws.Range("I55").Value = 1400
$ws.Range("J55").Value = 1550
$ws.Range("K55").Value = 1400
$ws.Range("L55").Value = 1550
$ws.Range("M55").Value = -1227
$ws.Range("N55").Value = -1896
$ws.Range("H68").Value = 6703.5454
$ws.Range("I68").Value = 6418
$ws.Range("J68").Value = 7465
$ws.Range("K68").Value = 6418
$ws.Range("L68").Value = 7465
$ws.Range("M68").Value = -5669
$ws.Range("N68").Value = -8963
$ws.Range("H71").Value = 6703.5454
$ws.Range("I71").Value = 6418
$ws.Range("J71").Value = 7465
$ws.Range("K71").Value = 32090
$ws.Range("L71").Value = 37325
$ws.Range("M71").Value = -28346
$ws.Range("N71").Value = -44813
$ws.Range("H93").Value = 1899.6666
$ws.Range("J93").Value = 1850
$ws.Range("L93").Value = 1850
$ws.Range("N93").Value = -4346
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()  # was -64937
$ws.Range("H122").Value = 21749.666
$ws.Range("J122").Value = 27500
$ws.Range("L122").Value = 82500
$ws.Range("N122").Value = -87400
$ws.Range("H136").Value = 41669130
$ws.Range("I136").Value = 2754.4
$ws.Range("J136").Value = 250001000
$ws.Range("K136").Value = 8263.200000000001
$ws.Range("L136").Value = 750003000
$ws.Range("M136").Value = -5713.200000000001
$ws.Range("N136").Value = -750008100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3667.3333
$ws.Range("I62").Value = 3667.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3667.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3043.3333
$ws.Range("N62").ClearContents()  # was -2448
$ws.Range("H65").Value = 3667.3333
$ws.Range("I65").Value = 3667.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 18336.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -15216.6665
$ws.Range("N65").ClearContents()  # was -12240
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
